# "Generate Report for Archive"
#
# The localization-status report is refreshed: the "Ready for handoff"
# status (shown for the zh-cn / de-de targets) has moved on to
# "In Translation", and the Status columns are re-fit to the new
# (shorter) text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: per-locale status columns (E = zh-cn, F = de-de) ---
if ($wsOverview.Range("E2").Text -eq $oldStatus) {
    $wsOverview.Range("E2").Value = $newStatus
}
if ($wsOverview.Range("F2").Text -eq $oldStatus) {
    $wsOverview.Range("F2").Value = $newStatus
}

# --- Per-locale detail sheets: Status column (C) ---
if ($wsZhCn.Range("C2").Text -eq $oldStatus) {
    $wsZhCn.Range("C2").Value = $newStatus
}
if ($wsDeDe.Range("C2").Text -eq $oldStatus) {
    $wsDeDe.Range("C2").Value = $newStatus
}

# --- Re-fit the Status columns now that the text is shorter ---
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
